$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "202.173.124.126"
$ws.Range("B9").Value = 28.3621642
$ws.Range("C9").Value = 77.2827832
$ws.Range("D9").Value = 13.041999816894531
$ws.Range("E9").Value = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
$ws.Range("F9").Value = "Linux armv81"
$ws.Range("G9").Value = "2025-06-25T16:29:33.084Z"
